# Updated CHE model - 2025-09-06 17:10
# Edits the "Misc" sheet of the SubRES_New_RE_and_Conventional_Trans workbook:
#  - ELC_won* -> ELC_wo* (now a stem matching both offshore & onshore units)
#  - E[_]W* / wind (TFM_AVA) row split into two explicit rows:
#      E[_]WOF* / windoff   (offshore)
#      E[_]WON* / windon    (onshore, newly inserted)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# Create the new shared strings in the same order they first appear in the
# target workbook ("windoff", "windon", then "ELC_wo*") so the rebuilt
# shared-string table lines up with the authored edit.

# Row 43: E[_]W* / wind  ->  E[_]WOF* / windoff
$ws.Range("C43").Value = "E[_]WOF*"
$ws.Range("D43").Value = "windoff"

# Insert a new row 44 (pushes the old rows 44+ down by one) and populate it
# with the new onshore-wind TFM_AVA entry: E[_]WON* / windon
$ws.Rows.Item(44).Insert()
$ws.Range("C44").Value = "E[_]WON*"
$ws.Range("D44").Value = "windon"
$ws.Range("E44").Value = "IN"

# Row 41: ELC_won* -> ELC_wo* (now matches both offshore and onshore units)
$ws.Range("D41").Value = "ELC_wo*"

# Keep the visible selection consistent with the authored workbook state
$ws.Range("D42").Select()
